$d = $word.ActiveDocument

# --- Title: April -> July ---
$null = $d.Content.Find.Execute("أفريل", $true, $false, $false, $false, $false, $true, 1, $false, "جويلية", 2)

# --- Table numeric updates ---

$t1 = $d.Tables.Item(1)
$t1.Rows.Item(2).Cells.Item(3).Range.Text = "931"
$t1.Rows.Item(2).Cells.Item(4).Range.Text = "9 310 000,00"
$t1.Rows.Item(2).Cells.Item(5).Range.Text = "10 980 000,00"
$t1.Rows.Item(3).Cells.Item(3).Range.Text = "167"
$t1.Rows.Item(3).Cells.Item(4).Range.Text = "1 670 000,00"
$t1.Rows.Item(4).Cells.Item(3).Range.Text = "205"
$t1.Rows.Item(4).Cells.Item(4).Range.Text = "2 050 000,00"
$t1.Rows.Item(4).Cells.Item(5).Range.Text = "5 070 000,00"
$t1.Rows.Item(5).Cells.Item(3).Range.Text = "160"
$t1.Rows.Item(5).Cells.Item(4).Range.Text = "1 600 000,00"
$t1.Rows.Item(6).Cells.Item(3).Range.Text = "94"
$t1.Rows.Item(6).Cells.Item(4).Range.Text = "940 000,00"
$t1.Rows.Item(7).Cells.Item(3).Range.Text = "48"
$t1.Rows.Item(7).Cells.Item(4).Range.Text = "480 000,00"
$t1.Rows.Item(8).Cells.Item(3).Range.Text = "208"
$t1.Rows.Item(8).Cells.Item(4).Range.Text = "2 080 000,00"
$t1.Rows.Item(8).Cells.Item(5).Range.Text = "5 090 000,00"
$t1.Rows.Item(9).Cells.Item(3).Range.Text = "127"
$t1.Rows.Item(9).Cells.Item(4).Range.Text = "1 270 000,00"
$t1.Rows.Item(10).Cells.Item(3).Range.Text = "68"
$t1.Rows.Item(10).Cells.Item(4).Range.Text = "680 000,00"
$t1.Rows.Item(11).Cells.Item(3).Range.Text = "46"
$t1.Rows.Item(11).Cells.Item(4).Range.Text = "460 000,00"
$t1.Rows.Item(12).Cells.Item(3).Range.Text = "60"
$t1.Rows.Item(12).Cells.Item(4).Range.Text = "600 000,00"
$t1.Rows.Item(13).Cells.Item(3).Range.Text = "373"
$t1.Rows.Item(13).Cells.Item(4).Range.Text = "3 730 000,00"
$t1.Rows.Item(13).Cells.Item(5).Range.Text = "4 980 000,00"
$t1.Rows.Item(14).Cells.Item(3).Range.Text = "33"
$t1.Rows.Item(14).Cells.Item(4).Range.Text = "330 000,00"
$t1.Rows.Item(16).Cells.Item(3).Range.Text = "48"
$t1.Rows.Item(16).Cells.Item(4).Range.Text = "480 000,00"
$t1.Rows.Item(17).Cells.Item(3).Range.Text = "167"
$t1.Rows.Item(17).Cells.Item(4).Range.Text = "1 670 000,00"
$t1.Rows.Item(17).Cells.Item(5).Range.Text = "4 170 000,00"
$t1.Rows.Item(18).Cells.Item(3).Range.Text = "27"
$t1.Rows.Item(18).Cells.Item(4).Range.Text = "270 000,00"
$t1.Rows.Item(19).Cells.Item(3).Range.Text = "121"
$t1.Rows.Item(19).Cells.Item(4).Range.Text = "1 210 000,00"
$t1.Rows.Item(20).Cells.Item(3).Range.Text = "102"
$t1.Rows.Item(20).Cells.Item(4).Range.Text = "1 020 000,00"
$t1.Rows.Item(21).Cells.Item(3).Range.Text = "95"
$t1.Rows.Item(21).Cells.Item(4).Range.Text = "950 000,00"
$t1.Rows.Item(21).Cells.Item(5).Range.Text = "3 870 000,00"
$t1.Rows.Item(22).Cells.Item(3).Range.Text = "137"
$t1.Rows.Item(22).Cells.Item(4).Range.Text = "1 370 000,00"
$t1.Rows.Item(23).Cells.Item(3).Range.Text = "98"
$t1.Rows.Item(23).Cells.Item(4).Range.Text = "980 000,00"
$t1.Rows.Item(25).Cells.Item(2).Range.Text = "3416"
$t1.Rows.Item(25).Cells.Item(3).Range.Text = "34 160 000,00"
$t1.Rows.Item(25).Cells.Item(4).Range.Text = "34 160 000,00"

$t2 = $d.Tables.Item(2)
$t2.Rows.Item(2).Cells.Item(3).Range.Text = "475"
$t2.Rows.Item(2).Cells.Item(4).Range.Text = "4 750 000,00"
$t2.Rows.Item(2).Cells.Item(5).Range.Text = "6 100 000,00"
$t2.Rows.Item(3).Cells.Item(3).Range.Text = "83"
$t2.Rows.Item(3).Cells.Item(4).Range.Text = "830 000,00"
$t2.Rows.Item(4).Cells.Item(3).Range.Text = "52"
$t2.Rows.Item(4).Cells.Item(4).Range.Text = "520 000,00"
$t2.Rows.Item(5).Cells.Item(3).Range.Text = "212"
$t2.Rows.Item(5).Cells.Item(4).Range.Text = "2 120 000,00"
$t2.Rows.Item(5).Cells.Item(5).Range.Text = "2 800 000,00"
$t2.Rows.Item(6).Cells.Item(3).Range.Text = "68"
$t2.Rows.Item(6).Cells.Item(4).Range.Text = "680 000,00"
$t2.Rows.Item(7).Cells.Item(2).Range.Text = "890"
$t2.Rows.Item(7).Cells.Item(3).Range.Text = "8 900 000,00"
$t2.Rows.Item(7).Cells.Item(4).Range.Text = "8 900 000,00"
$t2.Rows.Item(8).Cells.Item(2).Range.Text = "4306"
$t2.Rows.Item(8).Cells.Item(3).Range.Text = "43 060 000,00"
$t2.Rows.Item(8).Cells.Item(4).Range.Text = "43 060 000,00"

# --- Arabic words amount ---
$null = $d.Content.Find.Execute("أربعون مليون وخمسمئة وأربعون ألف", $true, $false, $false, $false, $false, $true, 1, $false, "ثلاثة وأربعون مليون وستون ألف", 2)
